# Weekly update: insert a new price record for "Hortaliza, Femacal de La
# Calera - Choclo" at row 665, shifting the existing rows 665-753 down to
# 666-754 (dimension grows from A1:R753 to A1:R754).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 665; Excel shifts 665:753 down to 666:754 and
# copies formatting (e.g. the date style on column D) from the row above.
$ws.Rows("665:665").Insert()

# Populate the new row 665 with the new weekly record.
$ws.Cells.Item(665, 1).Value  = 3
$ws.Cells.Item(665, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(665, 3).Value  = "Coquimbo"
$ws.Cells.Item(665, 4).Value  = 44776
$ws.Cells.Item(665, 5).Value  = 5
$ws.Cells.Item(665, 6).Value  = 100112024
$ws.Cells.Item(665, 7).Value  = "Choclo"
$ws.Cells.Item(665, 8).Value  = "Dulce o Americano"
$ws.Cells.Item(665, 9).Value  = "Primera"
$ws.Cells.Item(665, 10).Value = 85
$ws.Cells.Item(665, 11).Value = 40000
$ws.Cells.Item(665, 12).Value = 41000
$ws.Cells.Item(665, 13).Value = 40471
$ws.Cells.Item(665, 14).Value = "$/malla 70 unidades"
$ws.Cells.Item(665, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(665, 16).Value = 578
$ws.Cells.Item(665, 17).Value = 70
$ws.Cells.Item(665, 18).Value = "Hortaliza"

# Make sure column D keeps the workbook's date/time display format used by
# the rest of the "Fecha" column (row 666 is the row that used to be 665,
# so it still carries the original number format).
$ws.Cells.Item(665, 4).NumberFormat = $ws.Cells.Item(666, 4).NumberFormat
